$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 27: "Query Exemplars by Label" - status tbd -> x
$ws.Range("B27").Value = "x"

# Row 31: "Sort: exemplars by number of users" - add note
$ws.Range("C31").Value = "by ratings?"

# Row 32: "Sort: Top contributors acc. To label" - status tbd -> x
$ws.Range("B32").Value = "x"

# Row 33: "Sort: Contributors by ratings of contributed exemplars" - status tbd -> x
$ws.Range("B33").Value = "x"

# Update sheet view selection to C32 (matches final cursor position after edits)
$ws.Range("C32").Select()
